$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 6, column F (Classes) ---
$ws.Range("F6").Value = "4.GINFO"

# --- Add new rows 7, 8, 9 ---
$ws.Range("D7").Value = "Java"
$ws.Range("E7").Value = "profmail4@gmail.com"
$ws.Range("F7").Value = "3.GINFO"

$ws.Range("D8").Value = "Catia"
$ws.Range("E8").Value = "profmail5@gmail.com"
$ws.Range("F8").Value = "Cp 2"

$ws.Range("D9").Value = "Analyse 2"
$ws.Range("E9").Value = "profmail6@gmail.com"
$ws.Range("F9").Value = "Cp 1"

# --- Update rows 4-5, column F (Classes) last, so this string is appended last ---
$ws.Range("F4").Value = "4.GINFO/4.GTR"
$ws.Range("F5").Value = "4.GINFO/4.GTR"

# --- Hyperlinks for the new emails, matching the style of existing ones ---
$ws.Hyperlinks.Add($ws.Range("E7"), "mailto:profmail4@gmail.com")
$ws.Hyperlinks.Add($ws.Range("E8"), "mailto:profmail5@gmail.com")
$ws.Hyperlinks.Add($ws.Range("E9"), "mailto:profmail6@gmail.com")

# --- Match the style applied to other "Ensiegnant_Email" cells ---
$ws.Range("E7").Style = $ws.Range("E6").Style
$ws.Range("E8").Style = $ws.Range("E6").Style
$ws.Range("E9").Style = $ws.Range("E6").Style

# --- Update selection to F5 ---
$ws.Range("F5").Select()
